$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'45.346.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.30%  '

# Row 3
$ws.Range("D3").Value = "'2.377.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.35%  '

# Row 4
$ws.Range("E4").Value = '  +0.38%  '

# Row 5
$ws.Range("D5").Value = "'112.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.41%  '

# Row 6
$ws.Range("D6").Value = "'317.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.23%  '

# Row 7
$ws.Range("D7").Value = "'0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.51%  '

# Row 8
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("E9").Value = '  +5.31%  '

# Row 10
$ws.Range("D10").Value = "'42.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +11.17%  '

# Row 11
$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.02%  '

# Row 12
$ws.Range("D12").Value = "'8.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.34%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = "'0.109"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.00%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = "'1.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.10%  '

# Row 15
$ws.Range("D15").Value = "'15.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.93%  '

# Row 16
$ws.Range("D16").Value = "'2.741.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.63%  '

# Row 17
$ws.Range("D17").Value = "'2.374.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.12%  '

# Row 18
$ws.Range("D18").Value = "'45.270.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.15%  '

# Row 19
$ws.Range("D19").Value = "'7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.34%  '

# Row 20
$ws.Range("E20").Value = '  +3.93%  '

# Row 21
$ws.Range("D21").Value = "'13.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.41%  '

# Row 22
$ws.Range("D22").Value = "'74.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.23%  '

# Row 23
$ws.Range("D23").Value = "'3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.72%  '

# Row 24
$ws.Range("D24").Value = "'269.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.73%  '

# Row 25
$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.89%  '

# Row 26
$ws.Range("E26").Value = '  -0.68%  '

# Row 27
$ws.Range("D27").Value = "'11.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.43%  '

# Row 28
$ws.Range("D28").Value = "'7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.08%  '

# Row 30
$ws.Range("D30").Value = "'39.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.68%  '

# Row 31
$ws.Range("D31").Value = "'22.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.37%  '

# Row 32
$ws.Range("D32").Value = "'0.0948"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.33%  '

# Row 33
$ws.Range("D33").Value = "'170.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.85%  '

# Row 34
$ws.Range("D34").Value = "'2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.16%  '

# Row 35
$ws.Range("D35").Value = "'0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.96%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = "'0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.57%  '

# Row 37
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = "'4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.81%  '

# Row 38
$ws.Range("E38").Value = '  +14.03%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.0364"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.67%  '

# Row 40
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").Value = "'3.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.36%  '

# Row 41
$ws.Range("D41").Value = "'1.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.95%  '

# Row 42
$ws.Range("D42").Value = "'105.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.87%  '

# Row 43
$ws.Range("D43").Value = "'0.240"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.34%  '

# Row 44
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = "'71.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.63%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = "'13.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.90%  '

# Row 46
$ws.Range("E46").Value = '  +0.26%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'116.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.28%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = "'5.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.12%  '

# Row 49
$ws.Range("D49").Value = "'1.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +19.78%  '

# Row 50
$ws.Range("D50").Value = "'9.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.28%  '

# Row 51
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = "'0.223"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +18.51%  '
